$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log row (row 21) for the new test-mail entry that was
# processed during this sync run, mirroring the structure of row 20.

$row = 21

$ws.Cells.Item($row, 1).Value  = "Testmail #5: Wil je deze klant bellen?"
$ws.Cells.Item($row, 2).Value  = "Geachte afzender,`nDank voor uw e-mail. Om u beter van dienst te kunnen zijn, zouden we graag wat meer details ontvangen over welke klant we dienen te benaderen en waarvoor. Kunt u ons de naam van de klant en de reden voor het contact geven? Op die manier kunnen we dit efficiënt afhandelen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$ws.Cells.Item($row, 3).Value  = "Wil je deze klant bellen?"
$ws.Cells.Item($row, 4).Value  = "mailmind.test@zohomail.eu"
$ws.Cells.Item($row, 5).Value  = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($row, 6).Value  = "2025-07-29 21:37:31"
$ws.Cells.Item($row, 7).Value  = "Ja"
$ws.Cells.Item($row, 8).Value  = "Nee"
$ws.Cells.Item($row, 9).Value  = "Ja"
$ws.Cells.Item($row, 10).Value = "Nee"
